$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the two little "time (s) / accel (g)" tables used to chart the
# Base and Surface acceleration series.

# Base series (columns A:B) - column A/B already carry number formats from
# the sheet's column styles, so the header cells just inherit them.
$ws.Range("A2").Value = "time (s)"
$ws.Range("B2").Value = "accel (g)"

# Surface series (columns D:E) - brand new columns, so number formats have
# to be applied explicitly to line the headers up with the Base columns.
$ws.Range("D1").Value = "Surface"
$ws.Range("A1").Value = "Base"

$ws.Range("D2").Value = "time (s)"
$ws.Range("D2").NumberFormat = "0.00"

$ws.Range("E2").Value = "accel (g)"
$ws.Range("E2").NumberFormat = "0.000000"

# Leave the selection the way it ended up after building the data/chart
$ws.Range("AC1").Select() | Out-Null
